# Update cryptocurrency price/volume data per Thu Feb 15 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.339.06"
$ws.Range("E2").Value = "  +1.90%  "

$ws.Range("D3").Value = "2.798.65"
$ws.Range("E3").Value = "  +1.85%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'347.06"
$ws.Range("E5").Value = "  +4.33%  "

$ws.Range("D6").Value = "'116.64"
$ws.Range("E6").Value = "  +1.23%  "

$ws.Range("E7").Value = "  +4.18%  "

$ws.Range("E9").Value = "  +3.54%  "

$ws.Range("D10").Value = "'43.05"
$ws.Range("E10").Value = "  +3.94%  "

$ws.Range("D11").Value = "'0.0859"
$ws.Range("E11").Value = "  +3.70%  "

$ws.Range("D12").Value = "'20.13"
$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("E13").Value = "  +1.65%  "

$ws.Range("D14").Value = "'7.89"
$ws.Range("E14").Value = "  +3.05%  "

$ws.Range("D15").Value = "3.237.57"
$ws.Range("E15").Value = "  +1.83%  "

$ws.Range("D16").Value = "2.810.15"
$ws.Range("E16").Value = "  +2.07%  "

$ws.Range("D17").Value = "'0.896"
$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("D18").Value = "52.205.28"
$ws.Range("E18").Value = "  +1.92%  "

$ws.Range("D19").Value = "'3.24"
$ws.Range("E19").Value = "  +7.61%  "

$ws.Range("D20").Value = "'7.18"
$ws.Range("E20").Value = "  +4.51%  "

$ws.Range("D21").Value = "'13.46"
$ws.Range("E21").Value = "  -1.99%  "

$ws.Range("E22").Value = "  +2.15%  "

$ws.Range("D23").Value = "'70.24"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").Value = "'270.25"
$ws.Range("E24").Value = "  -3.47%  "

$ws.Range("D25").Value = "'2.76"
$ws.Range("E25").Value = "  +4.73%  "

$ws.Range("D26").Value = "'26.74"
$ws.Range("E26").Value = "  -0.85%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").Value = "'10.23"
$ws.Range("E28").Value = "  -1.26%  "

$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = "  +0.92%  "

$ws.Range("E30").Value = "  -0.36%  "

$ws.Range("D31").Value = "'35.07"
$ws.Range("E31").Value = "  -2.35%  "

$ws.Range("E32").Value = "  +0.47%  "

$ws.Range("E33").Value = "  +1.86%  "

$ws.Range("E34").Value = "  +24.28%  "

$ws.Range("D35").Value = "'0.0827"
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").Value = "'2.13"
$ws.Range("E36").Value = "  +0.56%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("E38").Value = "  -0.85%  "

$ws.Range("D39").Value = "'18.88"
$ws.Range("E39").Value = "  -3.61%  "

$ws.Range("D40").Value = "'3.23"
$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("D41").Value = "'2.68"
$ws.Range("E41").Value = "  +19.57%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'23.58"
$ws.Range("E42").Value = "  -0.92%  "

$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "'128.01"
$ws.Range("E43").Value = "  -1.16%  "

$ws.Range("E44").Value = "  +1.89%  "

$ws.Range("E45").Value = "  +0.90%  "

$ws.Range("D46").Value = "'3.36"
$ws.Range("E46").Value = "  -1.16%  "

$ws.Range("D47").Value = "2.070.15"
$ws.Range("E47").Value = "  -2.08%  "

$ws.Range("D48").Value = "'2.36"
$ws.Range("E48").Value = "  +2.87%  "

$ws.Range("D49").Value = "'0.976"
$ws.Range("E49").Value = "  +13.17%  "

$ws.Range("D50").Value = "'5.55"

$ws.Range("D51").Value = "'8.98"
$ws.Range("E51").Value = "  -0.86%  "
